$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 264 (shifts existing rows 264-315 down to 266-317)
$ws.Range("A264:A265").EntireRow.Insert()

# Fill the two new rows (264, 265) with copied constant metadata columns (A,B,C,E-J)
# copied from the surrounding rows, plus the new data values for this edit.
$newRows = @(264, 265)
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102005
    $ws.Cells.Item($r, 10).Value = "Naranja"
    $ws.Cells.Item($r, 20).Value = 15
}

# Row 264
$ws.Cells.Item(264, 4).Value = 44776
$ws.Cells.Item(264, 11).Value = "Fukumoto"
$ws.Cells.Item(264, 12).Value = "Primera"
$ws.Cells.Item(264, 13).Value = 200
$ws.Cells.Item(264, 14).Value = 6500
$ws.Cells.Item(264, 15).Value = 7000
$ws.Cells.Item(264, 16).Value = 6750
$ws.Cells.Item(264, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(264, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(264, 19).Value = 450

# Row 265
$ws.Cells.Item(265, 4).Value = 44776
$ws.Cells.Item(265, 11).Value = "Fukumoto"
$ws.Cells.Item(265, 12).Value = "Segunda"
$ws.Cells.Item(265, 13).Value = 100
$ws.Cells.Item(265, 14).Value = 5500
$ws.Cells.Item(265, 15).Value = 5500
$ws.Cells.Item(265, 16).Value = 5500
$ws.Cells.Item(265, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(265, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(265, 19).Value = 367
